$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row (row 1) shared-string labels:
#        "<Name>_old" -> "<Name>_FV2410"  (columns A:J)
#        "<Name>_new" -> "<Name>_FV2504"  (columns L:U)
#     Column K ("diff") is left untouched.
$headerBases = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $headerBases.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headerBases[$i] + "_FV2410"
}
for ($i = 0; $i -lt $headerBases.Length; $i++) {
    $ws.Cells.Item(1, 12 + $i).Value = $headerBases[$i] + "_FV2504"
}

# --- 2) Turn the used range A1:U54 into an Excel Table (ListObject) so the
#        header row gets filter buttons / structured-reference column names.
$rng = $ws.Range("A1:U54")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $rng, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"

# --- 3) Freeze the header row (split/freeze after row 1).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
